# Applies "Category App Setup Done" edit:
#  - colors every existing paragraph/run red (FF0000)
#  - prepends a new "Github / Initial branch / Website_setup" list block
#  - appends a new "Github / Adding drf  branch / Website_setup" list block
#
$d = $word.ActiveDocument
$W = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# ---------------------------------------------------------------------
# 1) Turn every existing paragraph (and its runs) red. Paragraphs(1..N-1)
#    happily accept Font.Color; the very last paragraph mark in the body
#    does not expose character formatting this way, so it is handled
#    below (step 3) via an XML replace instead.
# ---------------------------------------------------------------------
$total = $d.Paragraphs.Count
for ($i = 1; $i -lt $total; $i++) {
    $d.Paragraphs($i).Range.Font.Color = 255
}

# ---------------------------------------------------------------------
# 2) Prepend the new "Github" block before the very first paragraph.
# ---------------------------------------------------------------------
$startXml = @"
<w:p xmlns:w="$W"><w:pPr><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr></w:p>
<w:p xmlns:w="$W"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>Github</w:t></w:r></w:p>
<w:p xmlns:w="$W"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>Initial branch</w:t></w:r></w:p>
<w:p xmlns:w="$W"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>Website_setup</w:t></w:r></w:p>
<w:p xmlns:w="$W"><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:left="2160"/><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr></w:p>
"@
$d.Range(0, 0).InsertXML($startXml)

# ---------------------------------------------------------------------
# 3) Replace the original last (empty) paragraph in place with itself
#    (now red) plus the whole new closing block appended right after it
#    ("Github / Adding drf  branch / Website_setup" + trailing blanks).
#    Doing this as a single InsertXML on that paragraph's own Range
#    expands it in place, right before </w:body>/<w:sectPr>, which is
#    exactly the slot the new paragraphs belong in.
# ---------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$endXml = @"
<w:p xmlns:w="$W"><w:pPr><w:pStyle w:val="ListParagraph"/><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr></w:p>
<w:p xmlns:w="$W"><w:pPr><w:pStyle w:val="ListParagraph"/><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr></w:p>
<w:p xmlns:w="$W"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>Github</w:t></w:r></w:p>
<w:p xmlns:w="$W"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>Adding drf  branch</w:t></w:r></w:p>
<w:p xmlns:w="$W"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="2"/><w:numId w:val="1"/></w:numPr><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>Website_setup</w:t></w:r></w:p>
<w:p xmlns:w="$W"><w:pPr><w:pStyle w:val="ListParagraph"/><w:rPr><w:color w:val="FF0000"/></w:rPr></w:pPr></w:p>
"@
$lastPara.Range.InsertXML($endXml)

Write-Output "applied"
